$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data (A6, A7), matching style of A4/A5
$ws.Range("A6").Value = "8cdc5e926dfa8ac66d777503aaa28b7085e76b7c12e32dc7b7a578bc4892dea9"
$ws.Range("A7").Value = "3fbdafef1da1e7d10382fabfaec7d5980e9f6c765b953f5ba9e028d0c884f6d6"

$ws.Range("A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)

$ws.Range("A7").Select()
